$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.089.40'
$ws.Range('E2').Value = '  +3.47%  '
$ws.Range('D3').Value = '2.269.73'
$ws.Range('E3').Value = '  +1.52%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.33'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.636'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '63.41'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.93%  '
$ws.Range('E8').Value = '  +0.35%  '
$ws.Range('E9').Value = '  +10.34%  '
$ws.Range('E10').Value = '  +11.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.03'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '26.35'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +18.05%  '
$ws.Range('E13').Value = '  +2.11%  '
$ws.Range('D14').Value = '2.605.08'
$ws.Range('E14').Value = '  +1.68%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.64'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.22'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +9.59%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.843'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.35%  '
$ws.Range('D18').Value = '2.265.05'
$ws.Range('E18').Value = '  +1.35%  '
$ws.Range('D19').Value = '43.941.63'
$ws.Range('E20').Value = '  +7.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.63'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.08'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '252.71'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.06%  '
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.45'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.28'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.76%  '
$ws.Range('B27').Value = 'WEMIXToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.32'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +24.56%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.06'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '171.87'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.41%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.138'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.76%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.83'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.02%  '
$ws.Range('E32').Value = '  -5.47%  '
$ws.Range('E33').Value = '  +3.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0703'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.80'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.92'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.70%  '
$ws.Range('E37').Value = '  +6.81%  '
$ws.Range('E38').Value = '  +2.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.32'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0259'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.83%  '
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.70'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +9.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.000221'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.84%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0974'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.26'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.78%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '98.23'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.96%  '
$ws.Range('E47').Value = '  -0.55%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.24'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +11.35%  '
$ws.Range('E49').Value = '  -1.13%  '
$ws.Range('D50').Value = '1.447.46'
$ws.Range('E50').Value = '  -0.86%  '
$ws.Range('E51').Value = '  +3.87%  '
